$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns F (was E) through I (was H) one column to the
# right to make room for the new "Average insert length" column at E.
# Copy right-to-left so source cells aren't clobbered before they're read.
for ($row = 3; $row -le 15; $row++) {
    $srcH = $ws.Cells.Item($row, 8)   # old H
    $srcG = $ws.Cells.Item($row, 7)   # old G
    $srcF = $ws.Cells.Item($row, 6)   # old F
    $srcE = $ws.Cells.Item($row, 5)   # old E

    $dstI = $ws.Cells.Item($row, 9)   # new I
    $dstH = $ws.Cells.Item($row, 8)   # new H
    $dstG = $ws.Cells.Item($row, 7)   # new G
    $dstF = $ws.Cells.Item($row, 6)   # new F

    $dstI.Value2 = $srcH.Value2
    if ($row -gt 3) { $dstI.NumberFormat = $srcH.NumberFormat }

    $dstH.Value2 = $srcG.Value2
    if ($row -gt 3) { $dstH.NumberFormat = $srcG.NumberFormat }

    $dstG.Value2 = $srcF.Value2
    if ($row -gt 3) { $dstG.NumberFormat = $srcF.NumberFormat }

    $dstF.Value2 = $srcE.Value2
    if ($row -gt 3) { $dstF.NumberFormat = $srcE.NumberFormat }
}

# Header for new column E
$ws.Range("E3").Value2 = "Average insert length"

# New "Average insert length" data values for rows 4-15 (plain/General
# format - these cells inherited the old percentage style from column E's
# previous contents, so clear formatting back to the default first).
$values = @(2.1, 2.17, 2.1, 2.11, 2.25, 2.14, 2.09, 2.02, 2.08, 2.06, 2.13, 2.2)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 4 + $i
    $cell = $ws.Cells.Item($row, 5)
    $cell.ClearFormats()
    $cell.Value2 = $values[$i]
}

# Adjust column widths: H grows, and a new I width is set.
# (Input values compensate for this host's pixel-quantized ColumnWidth
# rounding so the stored width lands as close as possible to the target.)
$ws.Columns("H:H").ColumnWidth = 28.75
$ws.Columns("I:I").ColumnWidth = 26.33

# Update the selection to E16
$ws.Range("E16").Select()
